$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp
$ws.Range('A1').Value = 'Datos actualizados a 24 de Marzo de 2020 a las 13:16'

# Country data table (rows 4-199), sorted descending by Casos totales (col B)
# after refreshing stats for: Austria, Australia, Pakistan, Singapur, Rusia,
# Republica de Macedonia, Sri Lanka, Tunez
$data = @(
    @('China', 81171, 78, 73159, 4735, 1573, 7, 3277),
    @('Italia', 63927, 0, 7432, 50418, 3204, 0, 6077),
    @('Estados Unidos', 46168, 2434, 295, 45291, 1040, 29, 582),
    @('España', 39673, 4537, 3794, 33183, 2355, 385, 2696),
    @('Alemania', 30150, 1094, 453, 29567, 23, 7, 130),
    @('Iran', 24811, 1762, 8913, 13964, 0, 122, 1934),
    @('Francia', 19856, 0, 2200, 16796, 2082, 0, 860),
    @('Suiza', 9117, 322, 131, 8864, 141, 2, 122),
    @('Corea del Sur', 9037, 76, 3507, 5410, 59, 9, 120),
    @('Reino Unido', 6650, 0, 135, 6180, 20, 0, 335),
    @('Austria', 4809, 335, 9, 4775, 18, 4, 25),
    @('Paises Bajos', 4749, 0, 2, 4534, 435, 0, 213),
    @('Belgica', 4269, 526, 461, 3686, 381, 34, 122),
    @('Noruega', 2647, 22, 6, 2630, 42, 1, 11),
    @('Australia', 2144, 257, 118, 2018, 11, 1, 8),
    @('Canada', 2091, 0, 112, 1955, 1, 0, 24),
    @('Portugal', 2060, 0, 14, 2023, 47, 0, 23),
    @('Suecia', 2059, 13, 16, 2010, 110, 6, 33),
    @('Brasil', 1960, 36, 2, 1924, 18, 0, 34),
    @('Israel', 1656, 214, 49, 1606, 31, 0, 1),
    @('Malasia', 1624, 106, 183, 1426, 64, 1, 15),
    @('Dinamarca', 1577, 117, 1, 1544, 69, 8, 32),
    @('Turquia', 1529, 0, 0, 1492, 0, 0, 37),
    @('Chequia', 1289, 53, 7, 1280, 19, 1, 2),
    @('Japon', 1140, 12, 285, 813, 54, 0, 42),
    @('Irlanda', 1125, 0, 5, 1114, 29, 0, 6),
    @('Ecuador', 981, 0, 3, 960, 2, 0, 18),
    @('Pakistan', 916, 41, 13, 896, 0, 1, 7),
    @('Luxemburgo', 875, 0, 6, 861, 3, 0, 8),
    @('Tailandia', 827, 106, 52, 771, 7, 3, 4),
    @('Finlandia', 792, 92, 10, 781, 11, 0, 1),
    @('Polonia', 774, 25, 1, 764, 3, 1, 9),
    @('Rumania', 762, 186, 79, 675, 15, 1, 8),
    @('Chile', 746, 0, 11, 733, 7, 0, 2),
    @('Crucero', 712, 0, 587, 115, 15, 2, 10),
    @('Grecia', 695, 0, 29, 649, 35, 0, 17),
    @('Indonesia', 686, 107, 30, 601, 0, 6, 55),
    @('Islandia', 588, 0, 51, 535, 14, 1, 2),
    @('Arabia Saudita', 562, 0, 19, 543, 0, 0, 0),
    @('Singapur', 558, 49, 156, 400, 14, 0, 2),
    @('Sudafrica', 554, 152, 4, 550, 2, 0, 0),
    @('Filipinas', 552, 90, 20, 497, 1, 2, 35),
    @('India', 511, 12, 37, 464, 0, 0, 10),
    @('Catar', 501, 0, 37, 464, 6, 0, 0),
    @('Rusia', 495, 57, 22, 472, 8, 0, 1),
    @('Eslovenia', 442, 0, 0, 439, 12, 0, 3),
    @('Peru', 395, 0, 1, 389, 19, 0, 5),
    @('Barein', 390, 13, 164, 224, 3, 0, 2),
    @('Hong Kong', 386, 29, 102, 280, 4, 0, 4),
    @('Estonia', 369, 17, 7, 362, 4, 0, 0),
    @('Mexico', 367, 51, 4, 359, 1, 1, 4),
    @('Egipto', 366, 0, 68, 279, 0, 0, 19),
    @('Croacia', 361, 46, 5, 355, 6, 0, 1),
    @('Panama', 345, 0, 1, 338, 33, 0, 6),
    @('Colombia', 306, 29, 6, 297, 0, 0, 3),
    @('Libano', 304, 37, 8, 292, 4, 0, 4),
    @('Argentina', 301, 0, 51, 246, 0, 0, 4),
    @('Irak', 266, 0, 62, 181, 0, 0, 23),
    @('Serbia', 249, 0, 3, 243, 12, 0, 3),
    @('Republica Dominicana', 245, 0, 3, 239, 0, 0, 3),
    @('Armenia', 235, 0, 2, 233, 6, 0, 0),
    @('Argelia', 230, 0, 65, 148, 0, 0, 17),
    @('Taiwan', 216, 21, 29, 185, 0, 0, 2),
    @('Eslovaquia', 204, 18, 7, 197, 2, 0, 0),
    @('Bulgaria', 202, 1, 3, 196, 8, 0, 3),
    @('Emiratos Arabes Unidos', 198, 0, 41, 155, 2, 0, 2),
    @('Letonia', 197, 17, 1, 196, 0, 0, 0),
    @('Kuwait', 191, 2, 39, 152, 5, 0, 0),
    @('Lituania', 187, 8, 1, 185, 1, 0, 1),
    @('San Marino', 187, 0, 4, 163, 13, 0, 20),
    @('Hungria', 187, 20, 21, 157, 6, 1, 9),
    @('Principado de Andorra', 164, 31, 1, 162, 7, 0, 1),
    @('Uruguay', 162, 0, 0, 162, 3, 0, 0),
    @('Costa Rica', 158, 0, 2, 154, 2, 0, 2),
    @('Nueva Zelanda', 155, 53, 12, 143, 0, 0, 0),
    @('Bosnia y Herzegovina', 150, 14, 2, 147, 1, 0, 1),
    @('Republica de Macedonia', 148, 12, 1, 145, 1, 0, 2),
    @('Marruecos', 143, 0, 5, 134, 1, 0, 4),
    @('Jordania', 127, 0, 1, 126, 0, 0, 0),
    @('Albania', 123, 19, 10, 108, 2, 1, 5),
    @('Vietnam', 123, 0, 17, 106, 3, 0, 0),
    @('Islas Feroe', 122, 4, 23, 99, 0, 0, 0),
    @('Republica de Chipre', 116, 0, 3, 112, 3, 0, 1),
    @('Tunez', 114, 25, 1, 110, 11, 0, 3),
    @('Burkina Faso', 114, 15, 7, 103, 0, 0, 4),
    @('Malta', 110, 3, 2, 108, 1, 0, 0),
    @('Moldavia', 109, 0, 2, 106, 10, 0, 1),
    @('Brunei', 104, 13, 2, 102, 2, 0, 0),
    @('Sri Lanka', 101, 4, 2, 99, 2, 0, 0),
    @('Ucrania', 97, 24, 1, 93, 0, 0, 3),
    @('Camboya', 87, 0, 2, 85, 1, 0, 0),
    @('Senegal', 86, 7, 8, 78, 0, 0, 0),
    @('Venezuela', 84, 0, 15, 69, 2, 0, 0),
    @('Oman', 84, 18, 17, 67, 0, 0, 0),
    @('Bielorrusia', 81, 0, 22, 59, 0, 0, 0),
    @('Reunion', 75, 4, 1, 74, 0, 0, 0),
    @('Azerbaiyan', 72, 0, 10, 61, 3, 0, 1),
    @('Kazajistan', 68, 6, 0, 68, 0, 0, 0),
    @('Georgia', 67, 6, 9, 58, 1, 0, 0),
    @('Camerun', 66, 10, 2, 64, 0, 0, 0),
    @('Guadalupe', 62, 0, 0, 61, 4, 0, 1),
    @('Estado de Palestina', 60, 1, 17, 43, 0, 0, 0),
    @('Martinica', 53, 0, 0, 52, 7, 0, 1),
    @('Trinidad yTobago', 52, 1, 0, 52, 0, 0, 0),
    @('Ghana', 52, 25, 0, 50, 0, 0, 2),
    @('Liechtenstein', 51, 0, 0, 51, 0, 0, 0),
    @('Uzbekistan', 50, 4, 0, 50, 0, 0, 0),
    @('Consejo Danes para los Refugiados', 45, 9, 0, 43, 0, 0, 2),
    @('Kirguistan', 42, 26, 0, 42, 0, 0, 0),
    @('Afganistan', 42, 2, 1, 40, 0, 0, 1),
    @('Cuba', 40, 0, 0, 39, 3, 0, 1),
    @('Nigeria', 40, 0, 2, 37, 0, 0, 1),
    @('Puerto Rico', 39, 8, 1, 36, 0, 0, 2),
    @('Banglades', 39, 6, 5, 30, 0, 1, 4),
    @('Mayotte', 36, 12, 0, 36, 0, 0, 0),
    @('Ruanda', 36, 0, 0, 36, 0, 0, 0),
    @('Mauricio', 36, 0, 0, 34, 1, 0, 2),
    @('Honduras', 30, 0, 0, 30, 0, 0, 0),
    @('Guam', 29, 0, 0, 28, 0, 0, 1),
    @('Montenegro', 29, 2, 0, 28, 0, 0, 1),
    @('Bolivia', 28, 1, 0, 28, 0, 0, 0),
    @('Paraguay', 27, 5, 0, 25, 1, 1, 2),
    @('Costa de Marfil', 25, 0, 2, 23, 0, 0, 0),
    @('Macao', 25, 0, 10, 15, 0, 0, 0),
    @('Polinesia Francesa', 23, 5, 0, 23, 0, 0, 0),
    @('Monaco', 23, 0, 1, 22, 0, 0, 0),
    @('Guayana Francesa', 23, 3, 6, 17, 0, 0, 0),
    @('Jamaica', 21, 2, 2, 18, 0, 0, 1),
    @('Isla de Man', 20, 7, 0, 20, 0, 0, 0),
    @('Guatemala', 20, 0, 0, 19, 0, 0, 1),
    @('Guyana', 20, 0, 0, 19, 0, 0, 1),
    @('Togo', 18, 0, 0, 18, 0, 0, 0),
    @('Barbados', 17, 0, 0, 17, 0, 0, 0),
    @('Madagascar', 17, 5, 0, 17, 0, 0, 0),
    @('Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0),
    @('Kenia', 16, 0, 0, 16, 0, 0, 0),
    @('Gibraltar', 15, 0, 5, 10, 0, 0, 0),
    @('Maldivas', 13, 0, 5, 8, 0, 0, 0),
    @('Etiopia', 12, 1, 0, 12, 0, 0, 0),
    @('Tanzania', 12, 0, 0, 12, 0, 0, 0),
    @('Aruba', 12, 3, 1, 11, 0, 0, 0),
    @('Nueva Caledonia', 10, 2, 0, 10, 0, 0, 0),
    @('Mongolia', 10, 0, 0, 10, 0, 0, 0),
    @('Guinea Ecuatorial', 9, 0, 0, 9, 0, 0, 0),
    @('Uganda', 9, 0, 0, 9, 0, 0, 0),
    @('San Martin (Parte Francesa)', 8, 0, 0, 8, 0, 0, 0),
    @('Seychelles', 7, 0, 0, 7, 0, 0, 0),
    @('Bermudas', 6, 0, 0, 6, 0, 0, 0),
    @('Surinam', 6, 1, 0, 6, 0, 0, 0),
    @('Haiti', 6, 0, 0, 6, 0, 0, 0),
    @('Benin', 6, 0, 0, 6, 0, 0, 0),
    @('Gabon', 6, 0, 0, 5, 0, 0, 1),
    @('El Salvador', 5, 2, 0, 5, 0, 0, 0),
    @('Islas Caimanes', 5, 0, 0, 4, 0, 0, 1),
    @('Bahamas', 4, 0, 0, 4, 0, 0, 0),
    @('Suazilandia', 4, 0, 0, 4, 0, 0, 0),
    @('Guinea', 4, 0, 0, 4, 0, 0, 0),
    @('Fiyi', 4, 1, 0, 4, 0, 0, 0),
    @('Groenlandia', 4, 0, 0, 4, 0, 0, 0),
    @('Congo', 4, 0, 0, 4, 0, 0, 0),
    @('Namibia', 4, 0, 0, 4, 0, 0, 0),
    @('Curazao', 4, 0, 0, 3, 0, 0, 1),
    @('Republica del Chad', 3, 1, 0, 3, 0, 0, 0),
    @('Antigua y Barbuda', 3, 0, 0, 3, 0, 0, 0),
    @('Liberia', 3, 0, 0, 3, 0, 0, 0),
    @('Angola', 3, 0, 0, 3, 0, 0, 0),
    @('Niger', 3, 0, 0, 3, 0, 0, 0),
    @('Santa Lucia', 3, 0, 0, 3, 0, 0, 0),
    @('Zambia', 3, 0, 0, 3, 0, 0, 0),
    @('Republica de Africa Central', 3, 0, 0, 3, 0, 0, 0),
    @('Republica de Yibuti', 3, 0, 0, 3, 0, 0, 0),
    @('San Bartolome', 3, 0, 0, 3, 0, 0, 0),
    @('Cabo Verde', 3, 0, 0, 2, 0, 1, 1),
    @('Zimbabue', 3, 0, 0, 2, 0, 0, 1),
    @('Nicaragua', 2, 0, 0, 2, 0, 0, 0),
    @('Laos', 2, 2, 0, 2, 0, 0, 0),
    @('Butan', 2, 0, 0, 2, 0, 0, 0),
    @('Dominica', 2, 0, 0, 2, 0, 0, 0),
    @('Mauritania', 2, 0, 0, 2, 0, 0, 0),
    @('San Martin (Parte Holandesa)', 2, 0, 0, 2, 0, 0, 0),
    @('Birmania', 2, 0, 0, 2, 0, 0, 0),
    @('Sudan', 2, 0, 0, 1, 0, 0, 1),
    @('Nepal', 2, 0, 1, 1, 0, 0, 0),
    @('Gambia', 2, 0, 0, 1, 0, 0, 1),
    @('Montserrat', 1, 0, 0, 1, 0, 0, 0),
    @('San Vicente y las Granadinas', 1, 0, 0, 1, 0, 0, 0),
    @('Papua Nueva Guinea', 1, 0, 0, 1, 0, 0, 0),
    @('Timor Oriental', 1, 0, 0, 1, 0, 0, 0),
    @('Eritrea', 1, 0, 0, 1, 0, 0, 0),
    @('Mozambique', 1, 0, 0, 1, 0, 0, 0),
    @('Siria', 1, 0, 0, 1, 0, 0, 0),
    @('Granada', 1, 0, 0, 1, 0, 0, 0),
    @('Santa Sede', 1, 0, 0, 1, 0, 0, 0),
    @('Somalia', 1, 0, 0, 1, 0, 0, 0),
    @('Islas Turcas y Caicos', 1, 0, 0, 1, 0, 0, 0),
    @('Belice', 1, 0, 0, 1, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 4
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}